# Add csv module error handling
# Appends one new data row (row 30) to each of the four sheets, mirroring
# the existing "time / length / ID / ..." log layout.

$wb = $excel.ActiveWorkbook

$rowsData = @(
    @{
        Sheet = "MID_LFT_#1"
        A = 45816.46049768518
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
        D = "0x01,0x88"
        E = "0x07"
        F = 400
        G = [double]"5.68631262647113e+23"
        H = 392
        I = 7
    },
    @{
        Sheet = "MID_LFT_#2"
        A = 45816.46049768518
        B = "0x01,0x7c"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
        D = "0x01,0x78"
        E = "0x19"
        F = 380
        G = [double]"5.68432987514711e+23"
        H = 376
        I = 25
    },
    @{
        Sheet = "MID_PLT_#1"
        A = 45816.46049768518
        B = "0x00,0x6e"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
        D = "0x00,0x6D"
        E = "0x15"
        F = 110
        G = [double]"5.68631262647113e+23"
        H = 109
        I = 15
    },
    @{
        Sheet = "MID_PLT_#2"
        A = 45816.46049768518
        B = "0x00,0x82"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
        D = "0x00,0x81"
        E = "0x9"
        F = 130
        G = [double]"5.68631262647113e+23"
        H = 129
        I = 9
    }
)

foreach ($row in $rowsData) {
    $ws = $wb.Worksheets.Item($row.Sheet)
    $newRow = 30

    $ws.Cells.Item($newRow, 1).Value = $row.A
    $ws.Cells.Item($newRow, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($newRow, 2).Value = $row.B
    $ws.Cells.Item($newRow, 3).Value = $row.C
    $ws.Cells.Item($newRow, 4).Value = $row.D
    $ws.Cells.Item($newRow, 5).Value = $row.E

    $ws.Cells.Item($newRow, 6).Value = $row.F
    $ws.Cells.Item($newRow, 7).Value = $row.G
    $ws.Cells.Item($newRow, 8).Value = $row.H
    $ws.Cells.Item($newRow, 9).Value = $row.I
}
